# Add translation settings to ODK-X forms.
#
# The "settings" sheet gains five new columns (display.title.text.pt/.sw,
# display.locale.text, display.locale.text.pt/.sw) and three new rows that
# define the available locales (default/pt/sw -> English/Português/Kiswahili).
# The "settings" sheet also becomes the active sheet in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# --- Header row: new translation-related setting columns (D1:H1) ----------
$ws.Cells.Item(1, 4).Value = "display.title.text.pt"
$ws.Cells.Item(1, 5).Value = "display.title.text.sw"
$ws.Cells.Item(1, 6).Value = "display.locale.text"
$ws.Cells.Item(1, 7).Value = "display.locale.text.pt"
$ws.Cells.Item(1, 8).Value = "display.locale.text.sw"

# --- "survey" setting row: mirror the title into the new pt/sw columns ----
$ws.Cells.Item(5, 4).Value = "Water Bodies"
$ws.Cells.Item(5, 5).Value = "Water Bodies"

# --- New locale rows --------------------------------------------------------
# default locale -> English
$ws.Cells.Item(7, 1).Value = "default"
$ws.Cells.Item(7, 6).Value = "English"
$ws.Cells.Item(7, 7).Value = "English"
$ws.Cells.Item(7, 8).Value = "English"

# pt locale -> Português
$ws.Cells.Item(8, 1).Value = "pt"
$ws.Cells.Item(8, 6).Value = "Português"
$ws.Cells.Item(8, 7).Value = "Português"
$ws.Cells.Item(8, 8).Value = "Português"

# sw locale -> Kiswahili
$ws.Cells.Item(9, 1).Value = "sw"
$ws.Cells.Item(9, 6).Value = "Kiswahili"
$ws.Cells.Item(9, 7).Value = "Kiswahili"
$ws.Cells.Item(9, 8).Value = "Kiswahili"

# --- Make "settings" the active sheet / active cell -------------------------
$ws.Activate()
$ws.Range("A1").Select()
